$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Private")

# New row/cell values added to the "Private" sheet (matches the table that
# already lists WebMD / KKR / 2.8B for the Syapse deal with a new entry for
# Haiipfy's raise).
$ws.Range("B5").Value = "Haiipfy"
$ws.Range("E2").Value = "Raise"
$ws.Range("E5").Value = "73m"
$ws.Range("F2").Value = "Round"
$ws.Range("F5").Value = "D"

# The author switched focus to the "Private" sheet while making the edit.
$ws.Activate()
$ws.Range("F6").Select()
